# Q3 Update - 2025
# Applies the data refresh to the "fromCSV" sheet of the UN-MLI workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "short-url" placeholder in column B is repeated on every data row
# (rows 2-489). It was regenerated for this quarter's refresh.
$ws.Range("B2:B489").Value = "7qnKMH"

# Updated statistics (refugees / asylum_seekers / returned_refugees / idps /
# returned_idps) for a handful of rows. These columns hold numeric-looking
# text (not real numbers) in the source data, so force a text number format
# first -- otherwise Excel auto-converts the numeric-looking string into a
# real Number when it's assigned.
$statCells = @{
    "N470" = "93445"
    "N471" = "1385"
    "O471" = "117"
    "P471" = "16"
    "O472" = "22"
    "N473" = "18"
    "O473" = "7"
    "N475" = "301"
    "O475" = "50"
    "N478" = "15010"
    "Q479" = "360591"
    "R479" = "51546"
    "N480" = "25162"
    "O481" = "5"
    "O485" = "24"
    "O486" = "29"
    "N488" = "38"
    "O489" = "5"
}

# A1-untouched data cell that carries the same (unmodified) cell style as
# every target cell above, used below to restore formatting after writing.
$formatDonor = "N2"

foreach ($addr in $statCells.Keys) {
    # Force a text number format so the numeric-looking string is kept as
    # text (matching the source data's storage as text), then write it.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $statCells[$addr]
    # Re-copy the original (unchanged) cell formatting from an untouched
    # cell that shares the same style back onto this cell, so the edit
    # doesn't leave behind a stray number-format/style change.
    $ws.Range($formatDonor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
